# "added simple DSM - demand get be set up and down any time"
#
# 1) General Data sheet: annual PV capacity 20 -> 50, grid capacity limit 100 -> 20
# 2) irradiation sheet: replace the "Pv_Contractor"/"PV" profile (cols B/C) for hours
#    7-22 with a new profile (same new numbers as the "PVGIS_EU" series added below)
# 3) irradiation_winter sheet: add a new data series in column G, "PVGIS_EU", next to
#    the existing (now labeled) "PV_ninja_21_7" series in column F

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# General Data
# ---------------------------------------------------------------------------
$wsGen = $wb.Worksheets.Item("General Data")
$wsGen.Range("B3").Value = 50
$wsGen.Range("E3").Value = 20

# ---------------------------------------------------------------------------
# irradiation
# ---------------------------------------------------------------------------
$wsIrr = $wb.Worksheets.Item("irradiation")
$irrVals = @(0.1002,0.26207999999999998,0.41768,0.55237000000000003,0.48031999999999997,0.60196000000000005,0.48791000000000001,0.27662999999999999,0.19836000000000001,0.25089,0.21872999999999998,0.07665999999999999,0.0061600000000000005,0.0,0.0,0.0)
for ($i = 0; $i -lt $irrVals.Length; $i++) {
    $row = 7 + $i
    $wsIrr.Cells.Item($row, 2).Value = $irrVals[$i]
    $wsIrr.Cells.Item($row, 3).Value = $irrVals[$i]
}
$wsIrr.Range("D4").Select()

# ---------------------------------------------------------------------------
# irradiation_winter
# ---------------------------------------------------------------------------
$wsIrrWinter = $wb.Worksheets.Item("irradiation_winter")
$wsIrrWinter.Range("F1").Value = "PV_ninja_21_7"
$wsIrrWinter.Range("G1").Value = "PVGIS_EU"

$winterVals = @(0.0,0.0,0.0,0.0,0.0,0.1002,0.26207999999999998,0.41768,0.55237000000000003,0.48031999999999997,0.60196000000000005,0.48791000000000001,0.27662999999999999,0.19836000000000001,0.25089,0.21872999999999998,0.07665999999999999,0.0061600000000000005,0.0,0.0,0.0,0.0,0.0,0.0)
for ($i = 0; $i -lt $winterVals.Length; $i++) {
    $row = 2 + $i
    $wsIrrWinter.Cells.Item($row, 7).Value = $winterVals[$i]
}
$wsIrrWinter.Range("G2:G25").Select()
$excel.ActiveWindow.ScrollRow = 6

# ---------------------------------------------------------------------------
# Leave the originally active sheet/selection as the active one on save
# ---------------------------------------------------------------------------
$wsGen.Activate()
$wsGen.Range("B3").Select()
